$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# ---- Header row: rename/reindex the two "answer" columns (E1, F1) ----
$ws.Range("E1").Value = 'answer (shared_context)'
$ws.Range("F1").Value = 'answer (separate_context)'

# ---- The four Java snippets shown in column C (question content) ----
$findLastCode = '/**
    * Find last index of element
    * 
    *  @param x array to search
    *  @param y value to look for
    *  @return last index of y in x; -1 if absent
    *  @throws NullPointerException if x is null
    */
   public static int findLast (int[] x, int y)
   {       
    for (int i=x.length-1; i > 0; i--)
      {
         if (x[i] == y) 
         {
            return i;
         }
      }
      return -1;
   }'
$lastZeroCode = '/**
   * Find LAST index of zero
   *
   * @param x array to search
   * @return index of last 0 in x; -1 if absent
   * @throws NullPointerException if x is null
   */
   public static int lastZero (int[] x)
   {
      for (int i = 0; i < x.length; i++)
      {
         if (x[i] == 0)
         {
            return i;
         }
      }
      return -1;
   }'
$countPositiveCode = '/** 
   * Counts positive elements in array
   *
   * @param x array to search
   * @return number of positive elements in x
   * @throws NullPointerException if x is null
   */
   public static int countPositive (int[] x)
   {
      int count = 0;
   
      for (int i=0; i < x.length; i++)
      {
         if (x[i] >= 0)
         {
            count++;
         }
      }
      return count;
   }'
$oddOrPosCode = '/**
   * Count odd or positive elements in an array
   *
   * @param x array to search
   * @return count of odd or positive elements in x
   * @throws NullPointerException if x is null
   */
   public static int oddOrPos (int[] x)
   {  // Effects:  if x is null throw NullPointerException
      // else return the number of elements in x that
      //      are either odd or positive (or both)
      int count = 0;
   
      for (int i = 0; i < x.length; i++)
      {
         if (x[i]%2 == 1 || x[i] > 0)
         {
            count++;
         }
      }
      return count;
   }'

# ---- Newly-authored answers for the findLast (id 5.1) a/b rows ----
$findLastAnswer = 'There is a bug in the given code because the loop condition `i > 0` should be `i >= 0` in order to search all elements in the array `x`. The loop condition `i > 0` only searches the elements from index 1 to the last index, which means that the element at index 0 will not be searched. This can be fixed by modifying the loop condition to `i >= 0`. Here is the modified code:
```
public static int findLast (int[] x, int y)
{       
  for (int i=x.length-1; i >= 0; i--)
  {
     if (x[i] == y) 
     {
        return i;
     }
  }
  return -1;
}
```'
$testCaseAnswer = 'A test case that does not execute the fault would be to pass in an array `x` that does not contain the value `y`. For example:
```
int[] x = {1, 2, 3, 4, 5};
int y = 6;
assert(findLast(x, y) == -1);
```
This test case will not execute the fault because the value `y` is not present in the array `x`, so the loop will execute all the way to the end and return `-1` as expected.'

# ---- Column C previously only held the snippet on the "a" row of each group; ----
# ---- now every row (a-f) in a group repeats its group's snippet. ----
for ($r = 2; $r -le 7; $r++)   { $ws.Cells.Item($r, 3).Value = $findLastCode }
for ($r = 8; $r -le 13; $r++)  { $ws.Cells.Item($r, 3).Value = $lastZeroCode }
for ($r = 14; $r -le 19; $r++) { $ws.Cells.Item($r, 3).Value = $countPositiveCode }
for ($r = 20; $r -le 25; $r++) { $ws.Cells.Item($r, 3).Value = $oddOrPosCode }

# ---- New "separate_context" answers for rows 2 (a) and 3 (b) ----
$ws.Range("F2").Value = $findLastAnswer
$ws.Range("F3").Value = $testCaseAnswer

# ---- Wrap the long answer text and size the rows to fit it ----
$ws.Range("F2:F3").WrapText = $true
$ws.Rows.Item(2).RowHeight = 289
$ws.Rows.Item(3).RowHeight = 170

# ---- F1 should look like the rest of the header row (bold, bordered, no fill) ----
$ws.Range("A1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# ---- Column widths: column C loses its custom width, D/E/F get new widths ----
# (ColumnWidth is specified in character units; Excel stores it internally with a
#  small constant pixel-rounding offset of 5/6 of a character, which we subtract
#  here so the persisted <col width="..."> matches the target exactly.)
$widthOffset = 0.8333333333333333
$ws.Columns.Item(3).ColumnWidth = 8.83203125 - $widthOffset
$ws.Columns.Item(4).ColumnWidth = 93.83203125 - $widthOffset
$ws.Columns.Item(5).ColumnWidth = 75.5 - $widthOffset
$ws.Columns.Item(6).ColumnWidth = 106.5 - $widthOffset

# ---- Selection / scroll position ----
$ws.Range("F2").Select()

Write-Output "edit complete"
